# Add a new "R = 0.75" decay row (row 19) below the existing table, mirroring
# the structure of the other R-rows above it: column A holds the R value,
# column B the starting count (20000), column C the first decay step
# (=B*$A$), and columns D:AA a shared formula that repeatedly multiplies the
# previous column by the fixed R value in $A$19 - showing how many
# generations it takes for the number to fall below 1000.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# R value for this row and the starting population.
$ws.Range("A19").Value = 0.75
$ws.Range("B19").Value = 20000

# First step uses its own (non-shared) formula, like the other rows.
$ws.Range("C19").Formula = "=B19*`$A`$19"

# Remaining steps share one formula definition, copied across D19:AA19.
$ws.Range("D19:AA19").Formula = "=C19*`$A`$19"

# Match the italic styling used for this block (B19 general format,
# C19:AA19 integer format), reusing the workbook's existing style entries.
$ws.Range("B19:AA19").Font.Italic = $true
$ws.Range("C19:AA19").NumberFormat = "0"

# Leave the selection where the author left it after entering this row.
$ws.Range("L19").Select()
